$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = " net income (loss)"
$ws.Range("B2").Value = 862
$ws.Range("C2").Value = -775
$ws.Range("D2").Value = -1063

# Row 3
$ws.Range("A3").Value = " depreciation amortization and impairment"
$ws.Range("B3").Value = 2322
$ws.Range("C3").Value = 2154
$ws.Range("D3").Value = 1901

# Row 4
$ws.Range("A4").Value = " stock-based compensation"
$ws.Range("B4").Value = 1734
$ws.Range("C4").Value = 898
$ws.Range("D4").Value = 749

# Row 5
$ws.Range("A5").Value = " amortization of debt discounts and issuance costs"
$ws.Range("B5").Value = 180
$ws.Range("C5").Value = 188
$ws.Range("D5").Value = 159

# Row 6
$ws.Range("A6").Value = " inventory and purchase commitments write-downs"
$ws.Range("B6").Value = 202
$ws.Range("C6").Value = 193
$ws.Range("D6").Value = 85

# Row 7
$ws.Range("A7").Value = " loss on disposals of fixed assets"
$ws.Range("B7").Value = 117
$ws.Range("C7").Value = 146
$ws.Range("D7").Value = 162

# Row 8
$ws.Range("A8").Value = " foreign currency transaction net loss (gain)"
$ws.Range("B8").Value = 114
$ws.Range("C8").Value = -48
$ws.Range("D8").Value = -2

# Row 9
$ws.Range("A9").Value = " non-cash interest and other operating activities"
$ws.Range("B9").Value = 228
$ws.Range("C9").Value = 186
$ws.Range("D9").Value = 49

# Row 10
$ws.Range("A10").Value = " accounts receivable"
$ws.Range("B10").Value = -652
$ws.Range("C10").Value = -367
$ws.Range("D10").Value = -497

# Row 11
$ws.Range("A11").Value = " inventory"
$ws.Range("B11").Value = -422
$ws.Range("C11").Value = -429
$ws.Range("D11").Value = -1023

# Row 12
$ws.Range("A12").Value = " operating lease vehicles"
$ws.Range("B12").Value = -1072
$ws.Range("C12").Value = -764
$ws.Range("D12").Value = -215

# Row 13
$ws.Range("A13").Value = " prepaid expenses and other current assets"
$ws.Range("B13").Value = -251
$ws.Range("C13").Value = -288
$ws.Range("D13").Value = -82

# Row 14
$ws.Range("A14").Value = " other non-current assets."
$ws.Range("B14").Value = -344
$ws.Range("C14").Value = 115
$ws.Range("D14").Value = -207

# Row 15
$ws.Range("A15").Value = " accounts payable and accrued liabilities"
$ws.Range("B15").Value = 2102
$ws.Range("C15").Value = 646
$ws.Range("D15").Value = 1797

# Row 16
$ws.Range("A16").Value = " deferred revenue"
$ws.Range("B16").Value = 321
$ws.Range("C16").Value = 801
$ws.Range("D16").Value = 406

# Row 17
$ws.Range("A17").Value = " customer deposits"
$ws.Range("B17").Value = 7
$ws.Range("C17").Value = -58
$ws.Range("D17").Value = -96

# Row 18
$ws.Range("A18").Value = " other long-term liabilities"
$ws.Range("B18").Value = 495
$ws.Range("C18").Value = -5
$ws.Range("D18").Value = -25
